$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 59: add Category "Psi Focus" for Wild Talent feat
$ws.Range("B59").Value = "Psi Focus"

# Row 61: Improved Maneuver (Advanced Combat)
$ws.Range("A61").Value = "Improved Maneuver"
$ws.Range("E61").Value = "You do not provoke  attacks of opportunity when you use the *disarm, trip, sunder* or *overrun* powers."
$ws.Range("B61").Value = "Advanced Combat"

# Row 62: Flicking Disarm (Advanced Combat, prereq-ish column C -> Improved Maneuver)
$ws.Range("A62").Value = "Flicking Disarm"
$ws.Range("B62").Value = "Advanced Combat"
$ws.Range("C62").Value = "Improved Maneuver"
$ws.Range("E62").Value = "When you *disarm* a target, you can choose any unoccupied square adjacent to you or the target for the item to land in. "

# Row 63: Distant Maneuver (Advanced Combat, column C -> Improved Maneuver)
$ws.Range("A63").Value = "Distant Maneuver"
$ws.Range("B63").Value = "Advanced Combat"
$ws.Range("C63").Value = "Improved Maneuver"
$ws.Range("E63").Value = "You can *disarm, trip* or *sunder* using a ranged weapon."

# Update view to match final state
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 40
$ws.Range("E63").Select()
